# Generate Report for Handback
# Refresh the handback status report: a later handback/xliff-generate pass updated
# the "Priority" field and the handoff/handback/XLIFF-generation timestamps.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: column G = "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-09-02 10:20:24"
$wsOverview.Range("G5").Value = "2016-09-02 10:20:24"

# zh-cn sheet: column E = "Priority", H = "Correspond Handoff Datetime", K = "Correspond Handback DateTime"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-09-02 10:20:14"
$wsZhCn.Range("H5").Value = "2016-09-02 10:20:14"
$wsZhCn.Range("K2").Value = "2016-09-02 10:20:43"
$wsZhCn.Range("K5").Value = "2016-09-02 10:20:43"

# de-de sheet: column E = "Priority", H = "Correspond Handoff Datetime", K = "Correspond Handback DateTime"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-09-02 10:20:24"
$wsDeDe.Range("H5").Value = "2016-09-02 10:20:24"
$wsDeDe.Range("K2").Value = "2016-09-02 10:20:50"
$wsDeDe.Range("K5").Value = "2016-09-02 10:20:50"
